$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.886.23'
$ws.Range("E2").Value = '  -1.00%  '
$ws.Range("D3").Value = '1.879.29'
$ws.Range("E3").Value = '  -1.79%  '
$ws.Range("E4").Value = '  -0.29%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.88'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.91%  '
$ws.Range("E6").Value = '  -0.37%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4624'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.83%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3880'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.15%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07854'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.12%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9838'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.97%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.79'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.03%  '
$ws.Range("D12").Value = '1.909.51'
$ws.Range("E12").Value = '  -2.06%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.997'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.23%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.674'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.31%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06977'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.37%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '88.62'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.85%  '
$ws.Range("E17").Value = '  -0.35%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009935'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.63%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.93'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.002'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.34%  '
$ws.Range("D21").Value = '28.887.28'
$ws.Range("E21").Value = '  -1.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.269'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.61%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.49%  '
$ws.Range("E24").Value = '  +2.19%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '156.22'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.33%  '
$ws.Range("E26").Value = '  -1.90%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '5.908'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.67%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '117.65'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.74%  '
$ws.Range("E29").Value = '  -6.52%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09368'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.43%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.9029'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.43%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.273'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.88%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.319'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.42%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.250'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.93%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.05753'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.37%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.171'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.27%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02077'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.53%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.001'
$ws.Range("D38").Style = "Normal"
$ws.Range("E39").Value = '  -6.65%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5661'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.27%  '
$ws.Range("E41").Value = '  -2.63%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.692'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.69%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.94'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.67%  '
$ws.Range("E44").Value = '  -3.79%  '
$ws.Range("E45").Value = '  -2.51%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.07041'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.35%  '
$ws.Range("E47").Value = '  -3.98%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.543'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.70%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '112.47'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.62%  '
$ws.Range("E50").Value = '  -5.38%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '71.00'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.43%  '
